# Adds a new "LCD cursor" keyboard-shortcut row to the keyboardshortcuts sheet:
#   u  |  Toggle LCD cursor (off/profile/template)
# The new row is inserted right after the existing "z" / "Toggle xy cursor
# clamp mode ..." row (row 11), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 ("c" / "Shows/Hides Controls"),
# shifting it (and everything below) down to make room for the new entry.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the LCD cursor shortcut.
$ws.Cells.Item(11, 1).Value2 = "u"
$ws.Cells.Item(11, 2).Value2 = "Toggle LCD cursor (off/profile/template)"

# Match the slightly tighter row height used by its neighbouring rows.
$ws.Rows.Item(11).RowHeight = 13.8

# Reflect the new active cell/selection on this sheet.
$ws.Range("B11").Select()
